# Weekly fruit/vegetable price update: insert one new daily record.
# A new row of data (2022-07-25) is inserted after the current row 37,
# pushing the existing rows 38-62 down to 39-63 (dimension A1:R62 -> A1:R63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38 - shifts rows 38..62 down to 39..63
# and carries the date-format style from column D of the row above.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new record.
$ws.Range("A38").Value = 3
$ws.Range("B38").Value = "Femacal de La Calera"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = 44767
$ws.Range("E38").Value = 5
$ws.Range("F38").Value = 100112035
$ws.Range("G38").Value = "Bruselas (repollito)"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 45
$ws.Range("K38").Value = 15000
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = 15000
$ws.Range("N38").Value = "$/malla 15 kilos"
$ws.Range("O38").Value = "Provincia de Quillota"
$ws.Range("P38").Value = 1000
$ws.Range("Q38").Value = 15
$ws.Range("R38").Value = "Hortaliza"
